$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell-by-cell updates reflecting the refreshed crypto market snapshot.
# For Price (column D) values that are valid numeric literals, force the
# cell to Text format first so the original "inline string" semantics of
# the Price column (e.g. thousands-dot formatted strings) are preserved
# instead of Excel auto-converting them to numbers.

$ws.Range('D2').Value = '27.624.52'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').Value = '1.668.06'
$ws.Range('E3').Value = '  -3.23%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.46'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.18'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('E10').Value = '  -1.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0880'
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('D12').Value = '1.905.47'
$ws.Range('E12').Value = '  -3.18%  '
$ws.Range('D13').Value = '1.666.21'
$ws.Range('E14').Value = '  -3.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.561'
$ws.Range('E15').Value = '  -0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.67'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').Value = '27.616.66'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '242.63'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('E19').Value = '  -3.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.72'
$ws.Range('E20').Value = '  -4.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.51'
$ws.Range('E22').Value = '  -2.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.36'
$ws.Range('E23').Value = '  -3.05%  '
$ws.Range('E24').Value = '  -3.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.70'
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.23'
$ws.Range('E26').Value = '  -3.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.49'
$ws.Range('E27').Value = '  -1.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('E29').Value = '  -2.42%  '
$ws.Range('E30').Value = '  +3.30%  '
$ws.Range('E31').Value = '  -1.35%  '
$ws.Range('E32').Value = '  -2.30%  '
$ws.Range('D33').Value = '1.468.34'
$ws.Range('E33').Value = '  -2.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.13'
$ws.Range('E34').Value = '  -4.41%  '
$ws.Range('E35').Value = '  -5.04%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.37'
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.930'
$ws.Range('E37').Value = '  -2.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.578'
$ws.Range('E38').Value = '  -4.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0173'
$ws.Range('E39').Value = '  -1.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '69.72'
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('E41').Value = '  -4.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('E43').Value = '  -6.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.23'
$ws.Range('E44').Value = '  -2.80%  '
$ws.Range('D45').Value = '1.812.69'
$ws.Range('E45').Value = '  -3.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.789'
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.77'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.42'
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('E49').Value = '  -4.35%  '
$ws.Range('E50').Value = '  -1.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.92'
$ws.Range('E51').Value = '  -4.18%  '
